# The document has two logo pictures (BTec logo, Pearson logo), each
# inserted once into the "first page" header/footer and once into the
# "default" (primary) header/footer. Word's internal auto-numbering for
# these inline pictures' display names (wp:docPr/@name, mirrored by
# pic:cNvPr/@name) got swapped:
#   BTec logo   (header, .jpg): image1.jpg -> image2.jpg
#   Pearson logo (footer, .png): image2.png -> image1.png
#
# InlineShape objects don't expose a settable Name in the Word object
# model, so each picture is temporarily converted to a floating Shape
# (which does expose .Name), renamed, then converted back to an inline
# picture so the wp:inline layout is preserved.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Rename-LogoPicture($range, $newName) {
    $ishp = $range.InlineShapes.Item(1)
    $shp = $ishp.ConvertToShape()
    $shp.Name = $newName
    $shp.ConvertToInlineShape() | Out-Null
}

# Headers: BTec_Logo-Orange (image1.jpg -> image2.jpg)
Rename-LogoPicture $sec.Headers.Item(1).Range "image2.jpg"
Rename-LogoPicture $sec.Headers.Item(2).Range "image2.jpg"

# Footers: Pearson logo (image2.png -> image1.png)
Rename-LogoPicture $sec.Footers.Item(1).Range "image1.png"
Rename-LogoPicture $sec.Footers.Item(2).Range "image1.png"

Write-Output "Renamed 4 logo pictures."
